$d = $word.ActiveDocument
$section = $d.Sections(1)

# Footer "Primary" (word/footer1.xml) - Pearson logo, docPr id="3": image1.png -> image2.png
$footerPrimary = $section.Footers(1)
$pearsonShape1 = $footerPrimary.Range.InlineShapes(1)
$pearsonShape1.Name = "image2.png"

# Footer "FirstPage" (word/footer2.xml) - Pearson logo, docPr id="2": image1.png -> image2.png
$footerFirstPage = $section.Footers(2)
$pearsonShape2 = $footerFirstPage.Range.InlineShapes(1)
$pearsonShape2.Name = "image2.png"

# Header "FirstPage" (word/header1.xml) - BTec logo, docPr id="1": image2.jpg -> image1.jpg
$headerFirstPage = $section.Headers(2)
$btecShape = $headerFirstPage.Range.InlineShapes(1)
$btecShape.Name = "image1.jpg"
